# Update the "Förändrad" (Changed) date column C for rows 2-9
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
